# For each worksheet (one per paint color: "Default green", "Green",
# "Yellow", "Orange", "Brown", "Red", "Default Red", "Blue"), the
# "Values" column (B2:B11) was blank for all 10 cars. Fill each of
# those cells with the name of the sheet it lives on, for every
# worksheet in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $colorName = $ws.Name
    for ($row = 2; $row -le 11; $row++) {
        $ws.Cells.Item($row, 2).Value = $colorName
    }
}
